$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their text formatting so Excel does not
# auto-convert values like "295.68" or "1.00" into numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '40.129.79'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '2.208.82'
$ws.Range("E3").Value = '  -0.69%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '295.68'
$ws.Range("E5").Value = '  +1.31%  '
$ws.Range("D6").Value = '87.62'
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '0.471'
$ws.Range("E9").Value = '  -0.58%  '
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").Value = '52.10'
$ws.Range("E10").Value = '  +6.77%  '
$ws.Range("B11").Value = 'Avalanche'
$ws.Range("C11").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D11").Value = '30.83'
$ws.Range("E11").Value = '  +1.46%  '
$ws.Range("E12").Value = '  -0.08%  '
$ws.Range("E13").Value = '  +2.61%  '
$ws.Range("E14").Value = '  -1.59%  '
$ws.Range("D15").Value = '2.553.42'
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("D16").Value = '13.80'
$ws.Range("E16").Value = '  -1.04%  '
$ws.Range("D17").Value = '2.179.11'
$ws.Range("E17").Value = '  -2.03%  '
$ws.Range("E18").Value = '  +1.11%  '
$ws.Range("D19").Value = '40.053.80'
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("D20").Value = '0.0₃0886'
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").Value = '11.28'
$ws.Range("E21").Value = '  -0.76%  '
$ws.Range("E22").Value = '  -1.21%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = '235.35'
$ws.Range("E24").Value = '  -0.57%  '
$ws.Range("E26").Value = '  +0.48%  '
$ws.Range("E27").Value = '  -0.79%  '
$ws.Range("D28").Value = '23.18'
$ws.Range("E28").Value = '  +2.06%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '2.17'
$ws.Range("E29").Value = '  -0.53%  '
$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").Value = '9.32'
$ws.Range("E30").Value = '  +1.26%  '
$ws.Range("D31").Value = '156.39'
$ws.Range("E31").Value = '  -0.11%  '
$ws.Range("D32").Value = '32.03'
$ws.Range("E32").Value = '  +0.98%  '
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("E34").Value = '  +0.15%  '
$ws.Range("E35").Value = '  +4.29%  '
$ws.Range("E36").Value = '  -0.65%  '
$ws.Range("E37").Value = '  -0.73%  '
$ws.Range("E38").Value = '  +1.40%  '
$ws.Range("E39").Value = '  +2.97%  '
$ws.Range("D40").Value = '1.73'
$ws.Range("E40").Value = '  +2.13%  '
$ws.Range("D41").Value = '15.64'
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '2.075.05'
$ws.Range("E42").Value = '  -1.98%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '3.80'
$ws.Range("E43").Value = '  -1.51%  '
$ws.Range("D44").Value = '19.17'
$ws.Range("E44").Value = '  +7.36%  '
$ws.Range("E45").Value = '  +1.17%  '
$ws.Range("D46").Value = '9.96'
$ws.Range("E46").Value = '  +0.41%  '
$ws.Range("E47").Value = '  +5.47%  '
$ws.Range("D48").Value = '1.91'
$ws.Range("E48").Value = '  -10.89%  '
$ws.Range("D49").Value = '2.425.44'
$ws.Range("E49").Value = '  -0.40%  '
$ws.Range("D50").Value = '1.12'
$ws.Range("E50").Value = '  +1.68%  '
$ws.Range("E51").Value = '  +0.93%  '
